$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.312.93'
$ws.Range("E2").Value = '  -3.18%  '

$ws.Range("D3").Value = '2.256.61'
$ws.Range("E3").Value = '  -4.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '493.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.41%  '

$ws.Range("D9").Value = '2.281.96'
$ws.Range("E9").Value = '  -3.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0933'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.15%  '

$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.323'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.66'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.42%  '

$ws.Range("D14").Value = '2.676.31'
$ws.Range("E14").Value = '  -3.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.63%  '

$ws.Range("D16").Value = '54.255.26'
$ws.Range("E16").Value = '  -3.31%  '

$ws.Range("E17").Value = '  -3.04%  '

$ws.Range("D18").Value = '2.277.16'
$ws.Range("E18").Value = '  -3.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '299.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.48%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.02%  '

$ws.Range("E25").Value = '  +0.41%  '

$ws.Range("E26").Value = '  +1.18%  '

$ws.Range("D27").Value = '2.380.38'
$ws.Range("E27").Value = '  -4.06%  '

$ws.Range("E28").Value = '  +0.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '163.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.61%  '

$ws.Range("E31").Value = '  -2.45%  '

$ws.Range("D32").Value = '0.0₃0683'
$ws.Range("E32").Value = '  -3.58%  '

$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.83'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.30%  '

$ws.Range("E38").Value = '  +0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.856'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.99%  '

$ws.Range("E40").Value = '  -0.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.51'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.376'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.31%  '

$ws.Range("E43").Value = '  +1.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.93%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '126.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.47%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.47%  '

$ws.Range("E47").Value = '  -0.86%  '

$ws.Range("E48").Value = '  -2.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '239.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0480'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.04%  '

$ws.Range("E51").Value = '  -1.48%  '
